# Applies the cryptos-list price/volume refresh described in the commit.
# Values are written as literal text (matching the source t="inlineStr" cells),
# so each cell is forced to text format first to stop Excel from turning
# number-looking strings (e.g. "0.9998", "4.110") into floating point values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$updates = @(
    @('D2', '29.415.93'),
    @('D3', '1.849.77'),
    @('E3', '  +0.26%  '),
    @('D4', '0.9998'),
    @('E4', '  +0.14%  '),
    @('D5', '240.65'),
    @('E5', '  +0.16%  '),
    @('D6', '0.6303'),
    @('E6', '  -0.02%  '),
    @('E7', '  +0.05%  '),
    @('D8', '0.07697'),
    @('E8', '  +2.29%  '),
    @('D9', '0.2942'),
    @('E9', '  -0.42%  '),
    @('D10', '24.54'),
    @('E10', '  +0.56%  '),
    @('E11', '  +0.62%  '),
    @('D12', '1.851.52'),
    @('E12', '  +0.13%  '),
    @('D13', '5.027'),
    @('E13', '  +0.71%  '),
    @('D14', '0.00001087'),
    @('E14', '  +8.81%  '),
    @('D15', '0.6806'),
    @('E15', '  -0.27%  '),
    @('D16', '83.69'),
    @('E16', '  +1.07%  '),
    @('D17', '2.095.55'),
    @('E17', '  -0.43%  '),
    @('D18', '6.151'),
    @('E18', '  +0.62%  '),
    @('D19', '29.445.50'),
    @('E19', '  +0.23%  '),
    @('D20', '229.17'),
    @('E20', '  +0.61%  '),
    @('E21', '  +0.40%  '),
    @('D23', '7.452'),
    @('E24', '  +0.05%  '),
    @('D25', '157.41'),
    @('E25', '  +0.07%  '),
    @('D26', '0.1390'),
    @('D27', '8.356'),
    @('E27', '  +0.19%  '),
    @('E28', '  +0.16%  '),
    @('B29', 'PancakeSwap'),
    @('C29', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'),
    @('D29', '1.470'),
    @('E29', '  +0.56%  '),
    @('B30', 'Toncoin'),
    @('C30', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'),
    @('D30', '1.313'),
    @('E30', '  +4.83%  '),
    @('D31', '0.05718'),
    @('E31', '  +0.96%  '),
    @('D32', '4.110'),
    @('E32', '  -0.23%  '),
    @('D33', '4.051'),
    @('E33', '  +0.85%  '),
    @('E34', '  +0.67%  '),
    @('D35', '1.159'),
    @('E35', '  +0.40%  '),
    @('E36', '  -0.28%  '),
    @('E37', '  -0.19%  '),
    @('B38', 'MXToken'),
    @('C38', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'),
    @('D38', '2.778'),
    @('E38', '  -0.05%  '),
    @('B39', 'Maker'),
    @('C39', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'),
    @('D39', '1.230.61'),
    @('E39', '  -2.29%  '),
    @('E40', '  -0.43%  '),
    @('D41', '6.505'),
    @('D42', '0.9149'),
    @('E42', '  +0.68%  '),
    @('E43', '  +0.06%  '),
    @('B44', 'Quant'),
    @('C44', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'),
    @('D44', '101.51'),
    @('E44', '  +0.32%  '),
    @('B45', 'Aave'),
    @('C45', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'),
    @('D45', '66.32'),
    @('E45', '  +0.32%  '),
    @('B46', 'BabyDogeCoin'),
    @('C46', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'),
    @('D46', '0.00000000122'),
    @('E46', '  +3.42%  '),
    @('B47', 'Aptos'),
    @('C47', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'),
    @('D47', '7.160'),
    @('E47', '  +1.43%  '),
    @('B48', 'TheSandbox'),
    @('C48', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'),
    @('D48', '0.4015'),
    @('E48', '  -0.50%  '),
    @('B49', 'EnergySwap'),
    @('C49', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'),
    @('D49', '9.051'),
    @('E49', '  -0.65%  '),
    @('E50', '  +0.41%  '),
    @('B51', 'Algorand'),
    @('C51', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'),
    @('D51', '0.1124'),
    @('E51', '  -0.02%  ')
)

foreach ($u in $updates) {
    Set-TextValue $u[0] $u[1]
}
